$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.988.09"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "1.979.82"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.77"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("E6").Value = "  +2.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.06"
$ws.Range("E7").Value = "  +4.24%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.381"
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0797"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.90"
$ws.Range("E12").Value = "  +9.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.27"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.842"
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("D15").Value = "2.271.66"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("E16").Value = "  +3.52%  "
$ws.Range("D17").Value = "1.981.49"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").Value = "36.888.03"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.04"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "0.0₃0859"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.52"
$ws.Range("E24").Value = "  +2.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  +11.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.26"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.60"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("E30").Value = "  +17.07%  "
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.85"
$ws.Range("E32").Value = "  +3.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0621"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("E34").Value = "  +6.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.29"
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("E39").Value = "  -5.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0973"
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.94"
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.17"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0212"
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.59"
$ws.Range("E44").Value = "  +4.03%  "
$ws.Range("D45").Value = "1.369.81"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.97"
$ws.Range("E46").Value = "  +2.66%  "
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.22"
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.27"
$ws.Range("E50").Value = "  +6.21%  "
$ws.Range("E51").Value = "  +9.73%  "
